$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old "MARKS / A / B / Spacing / HC / ERRORS" row (row 5) down to row 6,
# opening up row 4 for the new "Wood Thickness:" field.
$ws.Rows("5").Insert()

# Build the new row 4 ("Wood Thickness:" label + "3x" value) by cloning the
# formatting already used for the Pattern row (row 3) so styles stay identical.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "Wood Thickness:"

$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = "3x"

# Pattern dropdown: default to "Staggered" and flip the list order.
$ws.Range("C3").Value = "Staggered"
$ws.Range("C3").Validation.Delete()
$ws.Range("C3").Validation.Add(3, 1, 1, '"Staggered,Non-Staggered"')

# New Wood Thickness dropdown list.
$ws.Range("C4").Validation.Delete()
$ws.Range("C4").Validation.Add(3, 1, 1, '"3x,2"""')

# Column B must now be wide enough to fit "Wood Thickness:" and becomes its
# own (non-merged) column width entry separate from column A.
$ws.Columns("B").ColumnWidth = 15.17

# Selection moves to C2.
$ws.Range("C2").Select() | Out-Null

Write-Output "Wood Nailer Briding and Bolts applied"
